$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "40.727.50"
$ws.Range("E2").Value = "  -3.10%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.130.37"
$ws.Range("E3").Value = "  -3.89%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.03"
$ws.Range("E5").Value = "  -3.66%  "

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.593"
$ws.Range("E6").Value = "  -5.57%  "

# Row 7: Solana
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.98"
$ws.Range("E7").Value = "  -7.05%  "

# Row 8: USDC
$ws.Range("E8").Value = "  +0.06%  "

# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.553"
$ws.Range("E9").Value = "  -8.82%  "

# Row 10: Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.33"
$ws.Range("E10").Value = "  -11.84%  "

# Row 11: Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0880"
$ws.Range("E11").Value = "  -7.78%  "

# Row 12: OKB
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.00"
$ws.Range("E12").Value = "  -7.82%  "

# Row 13: TRON
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0986"
$ws.Range("E13").Value = "  -4.90%  "

# Row 14: Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.45"
$ws.Range("E14").Value = "  -8.69%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.449.82"
$ws.Range("E15").Value = "  -3.93%  "

# Row 16: Chainlink
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.13"
$ws.Range("E16").Value = "  -1.08%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "2.150.44"
$ws.Range("E17").Value = "  -2.91%  "

# Row 18: Polygon
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.765"
$ws.Range("E18").Value = "  -8.56%  "

# Row 19: WrappedBTC
$ws.Range("D19").Value = "40.554.85"
$ws.Range("E19").Value = "  -3.22%  "

# Row 20: ShibaInu
$ws.Range("D20").Value = "0.0₃0968"
$ws.Range("E20").Value = "  -10.11%  "

# Row 21: Litecoin
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.68"
$ws.Range("E21").Value = "  -5.95%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.65"
$ws.Range("E22").Value = "  -9.16%  "

# Row 23: BitcoinCash
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "220.24"
$ws.Range("E23").Value = "  -4.42%  "

# Row 24: InternetComputer(DFINITY)
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.28"
$ws.Range("E24").Value = "  -11.86%  "

# Row 25: Dai
$ws.Range("E25").Value = "  -0.21%  "

# Row 26: ImmutableX
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.83"
$ws.Range("E26").Value = "  -12.16%  "

# Row 27: Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.35"
$ws.Range("E27").Value = "  -13.24%  "

# Row 28: WEMIXToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.16"
$ws.Range("E28").Value = "  -13.27%  "

# Row 29: Toncoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.13"
$ws.Range("E29").Value = "  -2.75%  "

# Row 30: PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.10"
$ws.Range("E30").Value = "  -7.57%  "

# Row 31: Monero
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.57"
$ws.Range("E31").Value = "  -0.95%  "

# Row 32: EthereumClassic
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.28"
$ws.Range("E32").Value = "  -5.99%  "

# Row 33: InjectiveProtocol
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.71"
$ws.Range("E33").Value = "  +0.08%  "

# Row 34: Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0733"
$ws.Range("E34").Value = "  -8.07%  "

# Row 35: Filecoin
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.92"
$ws.Range("E35").Value = "  -12.47%  "

# Row 36: Stellar
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.117"
$ws.Range("E36").Value = "  -6.34%  "

# Row 37: RenderToken
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.96"
$ws.Range("E37").Value = "  -7.49%  "

# Row 38: Kaspa
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0955"
$ws.Range("E38").Value = "  -12.95%  "

# Row 39: VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0271"
$ws.Range("E39").Value = "  -10.28%  "

# Row 40: LidoDAOToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.00"
$ws.Range("E40").Value = "  -5.89%  "

# Row 41: Celestia
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.16"
$ws.Range("E41").Value = "  -18.24%  "

# Row 42: THORChain
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.18"
$ws.Range("E42").Value = "  -8.44%  "

# Row 43: MultiversX
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "55.91"
$ws.Range("E43").Value = "  -15.41%  "

# Row 44: Algorand
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.180"
$ws.Range("E44").Value = "  -9.10%  "

# Row 45: FraxShare
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.05"
$ws.Range("E45").Value = "  -8.64%  "

# Row 46: Cronos
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0939"
$ws.Range("E46").Value = "  -6.62%  "

# Row 47: Aave
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "94.49"
$ws.Range("E47").Value = "  -10.21%  "

# Row 48: ARBITRUM
$ws.Range("E48").Value = "  -6.33%  "

# Row 49: TrustWalletToken
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.08"
$ws.Range("E49").Value = "  -7.34%  "

# Row 50: HuobiToken
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.59"
$ws.Range("E50").Value = "  -3.79%  "

# Row 51: RocketPoolETH
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.334.79"
$ws.Range("E51").Value = "  -3.76%  "
